$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'68.187.00"
$ws.Range("E2").Value = "  +0.41%  "
$ws.Range("D3").Value = "'3.257.61"
$ws.Range("E3").Value = "  -0.76%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "'582.91"
$ws.Range("E5").Value = "  +0.16%  "
$ws.Range("D6").Value = "'185.05"
$ws.Range("E6").Value = "  +1.06%  "
$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = "  -0.07%  "
$ws.Range("D8").Value = "'0.599"
$ws.Range("E8").Value = "  -0.70%  "
$ws.Range("E9").Value = "  -1.67%  "
$ws.Range("D10").Value = "'6.66"
$ws.Range("E10").Value = "  -1.01%  "
$ws.Range("D11").Value = "'0.419"
$ws.Range("E11").Value = "  +0.28%  "
$ws.Range("D12").Value = "'3.827.91"
$ws.Range("E12").Value = "  -0.66%  "
$ws.Range("E13").Value = "  -0.15%  "
$ws.Range("D14").Value = "'28.25"
$ws.Range("E14").Value = "  -1.93%  "
$ws.Range("D15").Value = "'68.235.72"
$ws.Range("E15").Value = "  +0.51%  "
$ws.Range("E16").Value = "  +0.70%  "
$ws.Range("D17").Value = "'3.284.95"
$ws.Range("E17").Value = "  +0.07%  "
$ws.Range("D18").Value = "'5.86"
$ws.Range("E18").Value = "  -0.05%  "
$ws.Range("D19").Value = "'13.63"
$ws.Range("E19").Value = "  +0.55%  "
$ws.Range("D20").Value = "'393.58"
$ws.Range("E20").Value = "  +4.25%  "
$ws.Range("D21").Value = "'7.70"
$ws.Range("E21").Value = "  +0.05%  "
$ws.Range("D22").Value = "'0.999"
$ws.Range("E22").Value = "  -0.16%  "
$ws.Range("D23").Value = "'71.28"
$ws.Range("E23").Value = "  -0.07%  "
$ws.Range("D24").Value = "'0.520"
$ws.Range("E24").Value = "  +1.01%  "
$ws.Range("E25").Value = "  -0.71%  "
$ws.Range("E26").Value = "  +4.44%  "
$ws.Range("D27").Value = "'9.81"
$ws.Range("E27").Value = "  -0.03%  "
$ws.Range("E28").Value = "  +0.04%  "
$ws.Range("E29").Value = "  -0.61%  "
$ws.Range("D30").Value = "'5.73"
$ws.Range("E30").Value = "  -0.02%  "
$ws.Range("D31").Value = "'22.94"
$ws.Range("E31").Value = "  -0.14%  "
$ws.Range("E32").Value = "  +2.69%  "
$ws.Range("E33").Value = "  +0.18%  "
$ws.Range("E35").Value = "  -1.97%  "
$ws.Range("D36").Value = "'162.80"
$ws.Range("E36").Value = "  +0.41%  "
$ws.Range("D37").Value = "'1.96"
$ws.Range("E37").Value = "  +6.09%  "
$ws.Range("D38").Value = "'0.824"
$ws.Range("E38").Value = "  -3.63%  "
$ws.Range("D39").Value = "'26.78"
$ws.Range("E39").Value = "  -1.20%  "
$ws.Range("E40").Value = "  -1.18%  "
$ws.Range("E41").Value = "  -3.64%  "
$ws.Range("D42").Value = "'2.50"
$ws.Range("E42").Value = "  -5.16%  "
$ws.Range("B43").Value = "Hedera"
$ws.Range("C43").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D43").Value = "'0.0690"
$ws.Range("E43").Value = "  +1.12%  "
$ws.Range("B44").Value = "Maker"
$ws.Range("C44").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D44").Value = "'2.656.44"
$ws.Range("E44").Value = "  -0.35%  "
$ws.Range("B45").Value = "InjectiveProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D45").Value = "'25.42"
$ws.Range("E45").Value = "  -1.77%  "
$ws.Range("B46").Value = "OKB"
$ws.Range("C46").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D46").Value = "'41.10"
$ws.Range("E46").Value = "  +0.26%  "
$ws.Range("D47").Value = "'339.85"
$ws.Range("E47").Value = "  -3.45%  "
$ws.Range("E48").Value = "  -0.97%  "
$ws.Range("D49").Value = "'6.36"
$ws.Range("E49").Value = "  +3.07%  "
$ws.Range("D50").Value = "'31.58"
$ws.Range("E50").Value = "  +1.26%  "
$ws.Range("D51").Value = "'0.991"
$ws.Range("E51").Value = "  -1.22%  "
